$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B19: convert from text "3" to a real numeric value 3
$ws.Range("B19").Value = 3

# Add new row 20 with annotation data
$ws.Range("A20").Value = "Ruilin"

# B20 must stay a text string "2" (not numeric), matching the source format
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "2"
$ws.Range("B20").Style = "Normal"

$ws.Range("C20").Value = "No clear novelty"
$ws.Range("D20").Value = "DFT"
$ws.Range("E20").Value = "OTH"
$ws.Range("F20").Value = "4efacd8b-a5d8-471d-9660-f5eb687b96fc"
$ws.Range("G20").Value = "Byni8NLHf_annotated.xlsx"
$ws.Range("H20").Value = "No clear novelty"
